# Generate Report for Archive
# Updates the localization status workbook: the two files
# "1b327bbb-a660-4809-8fbe-cea71ac709f4.md" and
# "2b3654f8-0614-4937-bd50-b6b3041dc000.md" have moved from
# "Ready for handoff" to "In Translation" on the Overview sheet as
# well as on each per-language detail sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
